# "Add golden chaas, masala chaas and improvised veg upma format"
#
# Sheet1 holds the "Base Recipe" adjustment table (rows 2-11). The
# servings count (column E) moves from 2 to 3 for every ingredient row,
# a couple of quantities/units get corrected, and the "curry leaves
# paster" row's quantity format switches to a centered whole-number
# style (it is now counted in discrete units, not measured out).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Recipe now serves 3 instead of 2 -> bump every row's servings figure.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Value = 3
}

# "ginger grated" (row 5): 4 gm -> 0.5 gm per serving.
$ws.Range("F5").Value = 0.5

# "curry leaves paster" (row 9): unit corrected from "gram"/" springs" wording
# to proper units, and quantity bumped from 2 to 12.
$ws.Range("C8").Value = "tsp"
$ws.Range("C9").Value = "count"
$ws.Range("F9").Value = 12

# Row 9's result is now a whole number of items, shown centered with no
# decimals.
$ws.Range("B9").NumberFormat = "0"
$ws.Range("B9").HorizontalAlignment = -4108

# Leave the selection where the author left off editing.
[void]$ws.Range("F6").Select()
